# Applies the "Risolto bug in GestioneOperaio per label e NuovoCostoOperaio" edit:
#   1. Paragraph "Verificare funzionamento..."  -> green-highlight the whole
#      paragraph (incl. the paragraph mark) and split " Frut Loc"" into
#      separate spell-checked runs ("Frut" / "Loc" each wrapped in their own
#      proofErr spellStart/spellEnd, like "Ord" already was). The stray
#      _GoBack bookmark that lived at the end of this paragraph is removed.
#   2. Paragraph "Aggiungere alla "Stampa Ord Frut Cant Excel"..." -> also
#      green-highlighted (incl. paragraph mark) and now carries the _GoBack
#      bookmark (re-inserted at its start).

$d = $word.ActiveDocument

$hl = '<w:rPr><w:highlight w:val="green"/></w:rPr>'
$wOpenXmlHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$wOpenXmlFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# Locate the two target paragraphs by their (unique) text.
$paraVerificare = $null
$paraAggiungere = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Verificare funzionamento*Stampa*") {
        $paraVerificare = $p
    }
    if ($p.Range.Text -like "Aggiungere alla*Stampa Ord Frut Cant Excel*") {
        $paraAggiungere = $p
    }
}

# ---------------------------------------------------------------------
# 1) "Verificare funzionamento ..." paragraph
# ---------------------------------------------------------------------
$pPrVerificare = '<w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="23"/></w:numPr><w:spacing w:line="256" w:lineRule="auto"/>' + $hl + '</w:pPr>'

$runsVerificare = ''
$runsVerificare += '<w:r>' + $hl + '<w:t xml:space="preserve">Verificare funzionamento &#8220;Stampa </w:t></w:r>'
$runsVerificare += '<w:proofErr w:type="spellStart"/>'
$runsVerificare += '<w:r>' + $hl + '<w:t>Ord</w:t></w:r>'
$runsVerificare += '<w:proofErr w:type="spellEnd"/>'
$runsVerificare += '<w:r>' + $hl + '<w:t xml:space="preserve"> </w:t></w:r>'
$runsVerificare += '<w:proofErr w:type="spellStart"/>'
$runsVerificare += '<w:r>' + $hl + '<w:t>Frut</w:t></w:r>'
$runsVerificare += '<w:proofErr w:type="spellEnd"/>'
$runsVerificare += '<w:r>' + $hl + '<w:t xml:space="preserve"> </w:t></w:r>'
$runsVerificare += '<w:proofErr w:type="spellStart"/>'
$runsVerificare += '<w:r>' + $hl + '<w:t>Loc</w:t></w:r>'
$runsVerificare += '<w:proofErr w:type="spellEnd"/>'
$runsVerificare += '<w:r>' + $hl + '<w:t>&#8221;</w:t></w:r>'

$pXmlVerificare = '<w:p w:rsidR="00E90BFB" w:rsidRDefault="00E90BFB" w:rsidP="002937FF">' + $pPrVerificare + $runsVerificare + '</w:p>'
$xmlVerificare = $wOpenXmlHeader + $pXmlVerificare + $wOpenXmlFooter

# The _GoBack bookmark sits right at the end of this paragraph's range; drop
# it first so InsertXML (which replaces the whole paragraph, mark included)
# doesn't have to straddle it.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$paraVerificare.Range.InsertXML($xmlVerificare)

# ---------------------------------------------------------------------
# 2) "Aggiungere alla ... Stampa Ord Frut Cant Excel ..." paragraph
# ---------------------------------------------------------------------
$pPrAggiungere = '<w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="23"/></w:numPr><w:spacing w:line="256" w:lineRule="auto"/>' + $hl + '</w:pPr>'
$bookmarkAggiungere = '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>'
$runAggiungere = '<w:r>' + $hl + '<w:t>Aggiungere alla &#8220;Stampa Ord Frut Cant Excel&#8221; la lista dei frutti (Non appartenenti ad un gruppo)</w:t></w:r>'

$pXmlAggiungere = '<w:p w:rsidR="00E90BFB" w:rsidRDefault="00E90BFB" w:rsidP="000E0AED">' + $pPrAggiungere + $bookmarkAggiungere + $runAggiungere + '</w:p>'
$xmlAggiungere = $wOpenXmlHeader + $pXmlAggiungere + $wOpenXmlFooter

$paraAggiungere.Range.InsertXML($xmlAggiungere)

Write-Output "Edit complete."
